$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# NOTE: the source text uses non-breaking spaces (U+00A0) around the budget
# figures, e.g. "100<nbsp>000<nbsp>$.". PowerShell's `+` concatenation here
# coerces a lone [char] operand to its numeric string form when joined with a
# non-empty string literal, so nbsp strings are built with [string]::Concat
# instead, which preserves the character correctly.

$replacements = @(
    @{ Old = [string]::Concat("000", $nbsp, "`$."); New = [string]::Concat("000", $nbsp, "USD.") },
    @{ Old = "Aquí tiene 10 posibles lemas"; New = "Aquí tienes 10 posibles lemas" },
    @{ Old = "Té Chai: La especia de la vida"; New = "Té chai: el sabor de la vida" },
    @{ Old = "Té Chai: un mundo de sabor en una taza"; New = "Té chai: un mundo de sabor en una taza" },
    @{ Old = "Té Chai: Descubrir la magia de la India"; New = "Té chai: descubre la magia de la India" },
    @{ Old = "Té Chai: La mezcla perfecta de salud y placer"; New = "Té chai: la mezcla perfecta de salud y placer" },
    @{ Old = "Té Chai: Más que solo té, una forma de vida"; New = "Té chai: más que té, una forma de vida" },
    @{ Old = "Té Chai: Una bebida para todas las estaciones y razones"; New = "Té chai: una bebida para todas las estaciones y razones" },
    @{ Old = "Té Chai: la máxima indulgencia para sus sentidos"; New = "Té chai: la máxima indulgencia para tus sentidos" },
    @{ Old = "Té Chai: Un dulce escape del día a día"; New = "Té chai: Un dulce evasión del día a día" },
    @{ Old = "Té Chai: Compartir la calidez, compartir el amor"; New = "Té chai: comparte la calidez, comparte el amor" },
    @{ Old = "Té Chai: Tratate con algo especial"; New = "Té chai: date un gusto con algo especial" }
)

foreach ($rep in $replacements) {
    $found = $d.Content.Find.Execute($rep.Old, $true, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2)
    if (-not $found) {
        Write-Output "WARNING: replacement not found for: $($rep.Old)"
    }
}
